$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "assistant to the device using.",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "assistant to the system.",
    2
)
